$d = $word.ActiveDocument

# Remove the w14:ligatures (historicalDiscretional) compatibility element
# from every run's rPr in the main document story, the footnote story and
# the endnote story. The Word object model has no direct property for the
# w14:ligatures extension, so each story's content is rewritten with
# InsertXML using the same markup minus the <w14:ligatures/> element.

$mainXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:cs="Helvetica" w:eastAsia="Helvetica" w:hAnsi="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Footnote</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:cs="Helvetica" w:eastAsia="Helvetica" w:hAnsi="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:vertAlign w:val="superscript"/></w:rPr><w:footnoteReference w:id="2"/></w:r><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:cs="Helvetica" w:eastAsia="Helvetica" w:hAnsi="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> and Endnote</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:cs="Helvetica" w:eastAsia="Helvetica" w:hAnsi="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:vertAlign w:val="superscript"/></w:rPr><w:endnoteReference w:id="2"/></w:r><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:cs="Helvetica" w:eastAsia="Helvetica" w:hAnsi="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> Test</w:t></w:r></w:p>
'@
$d.Content.InsertXML($mainXml)

$footnoteXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:tabs><w:tab w:val="start" w:pos="0"/><w:tab w:val="start" w:pos="400"/></w:tabs></w:pPr><w:r><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:cs="Helvetica" w:eastAsia="Helvetica" w:hAnsi="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:vertAlign w:val="superscript"/></w:rPr><w:footnoteRef/></w:r><w:r><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:cs="Helvetica" w:eastAsia="Helvetica" w:hAnsi="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">This is the footnote content.</w:t></w:r></w:p>
'@
$d.Footnotes(1).Range.InsertXML($footnoteXml)

$endnoteXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:tabs><w:tab w:val="start" w:pos="0"/><w:tab w:val="start" w:pos="400"/></w:tabs></w:pPr><w:r><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:cs="Helvetica" w:eastAsia="Helvetica" w:hAnsi="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:vertAlign w:val="superscript"/></w:rPr><w:endnoteRef/></w:r><w:r><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:cs="Helvetica" w:eastAsia="Helvetica" w:hAnsi="Helvetica"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">This is the endnote content.</w:t></w:r></w:p>
'@
$d.Endnotes(1).Range.InsertXML($endnoteXml)
